$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L ("Rejeito"), shifting the existing
# "Papelão e jornal" ... columns one place to the right.
$ws.Range("L1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("L1").Value = "Rejeito"

# Row 2: rewrite the "Reciclaveis" total formula to exclude the new
# "Rejeito" column, and give the new column its own formula.
$ws.Range("K2").Formula = "=SUM(AK2-J2-L2)"
$ws.Range("L2").Formula = "=SUM(Z2,Y2,AC2)"

# Rows 3-8 share the same two formulas (mirrors the original shared
# formulas that covered K3:K8).
$ws.Range("K3:K8").Formula = "=SUM(AK3-J3-L3)"
$ws.Range("L3:L8").Formula = "=SUM(Z3,Y3,AC3)"

# Restore the selection the author left the sheet on.
$ws.Range("K12").Select()
